$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeAPI")

# Insert a new column before column G: Authenticator (F) | [NEW: AuthenticationParameter] | RequestParameters (was G, now H) | ...
$ws.Columns.Item(7).Insert()

# --- New header cell G1 = "AuthenticationParameter" (style matches other header cells: bold, green fill, bordered, centered) ---
$g1 = $ws.Cells.Item(1, 7)
$g1.Value = "AuthenticationParameter"
$g1.Font.Bold = $true
$g1.Interior.Color = 5296274
$g1.Borders.LineStyle = 1
$g1.HorizontalAlignment = -4108
$g1.VerticalAlignment = -4108

# Column width for the newly inserted column (displayed width 24)
$ws.Columns.Item(7).ColumnWidth = 23.166666666666664

# --- Row 2 (TC001_GET_F1_Driver): populate new Authenticator-parameter cells F2, G2 ---
$f2 = $ws.Cells.Item(2, 6)
$f2.Value = "PreemptiveBasicAuthScheme"
$f2.Borders.LineStyle = 1
$f2.HorizontalAlignment = -4131
$f2.VerticalAlignment = -4160
$f2.WrapText = $true

$g2 = $ws.Cells.Item(2, 7)
$g2.Value = "Username_ToolsQA:Password_TestPassword"
$g2.Borders.LineStyle = 1
$g2.HorizontalAlignment = -4131
$g2.VerticalAlignment = -4160
$g2.WrapText = $true

# Update the active selection to G2 as in the diff
$ws.Range("G2").Select()
